$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 442  # F4: 438 -> 442
$ws.Cells.Item(5, 6).Value = 1813  # F5: 1809 -> 1813
$ws.Cells.Item(7, 6).Value = 2895  # F7: 2891 -> 2895
$ws.Cells.Item(9, 6).Value = 755  # F9: 754 -> 755
$ws.Cells.Item(10, 6).Value = 7278  # F10: 0 -> 7278
$ws.Cells.Item(11, 6).Value = 178  # F11: 177 -> 178
$ws.Cells.Item(13, 6).Value = 191  # F13: 190 -> 191
$ws.Cells.Item(14, 6).Value = 1660  # F14: 1655 -> 1660
$ws.Cells.Item(15, 6).Value = 1417  # F15: 1412 -> 1417
$ws.Cells.Item(16, 6).Value = 1261  # F16: 1258 -> 1261
$ws.Cells.Item(17, 6).Value = 127  # F17: 126 -> 127
$ws.Cells.Item(18, 6).Value = 127  # F18: 126 -> 127
$ws.Cells.Item(19, 6).Value = 3316  # F19: 3297 -> 3316
$ws.Cells.Item(20, 6).Value = 5771  # F20: 5764 -> 5771
$ws.Cells.Item(21, 6).Value = 5771  # F21: 5764 -> 5771
$ws.Cells.Item(22, 6).Value = 564  # F22: 560 -> 564
$ws.Cells.Item(23, 6).Value = 928  # F23: 923 -> 928
$ws.Cells.Item(26, 6).Value = 5736  # F26: 5721 -> 5736
$ws.Cells.Item(27, 6).Value = 319  # F27: 318 -> 319
$ws.Cells.Item(29, 6).Value = 3948  # F29: 3941 -> 3948
$ws.Cells.Item(30, 6).Value = 204  # F30: 199 -> 204
$ws.Cells.Item(31, 6).Value = 659  # F31: 657 -> 659
$ws.Cells.Item(32, 6).Value = 1809  # F32: 1806 -> 1809
$ws.Cells.Item(34, 6).Value = 242  # F34: 240 -> 242
$ws.Cells.Item(35, 6).Value = 13  # F35: 12 -> 13
$ws.Cells.Item(36, 6).Value = 142  # F36: 141 -> 142
$ws.Cells.Item(37, 6).Value = 72  # F37: 66 -> 72
$ws.Cells.Item(38, 6).Value = 300  # F38: 298 -> 300
$ws.Cells.Item(39, 6).Value = 1097  # F39: 1096 -> 1097
$ws.Cells.Item(41, 6).Value = 1807  # F41: 1801 -> 1807
$ws.Cells.Item(42, 6).Value = 69  # F42: 68 -> 69
$ws.Cells.Item(43, 6).Value = 342  # F43: 337 -> 342
$ws.Cells.Item(45, 6).Value = 989  # F45: 985 -> 989
$ws.Cells.Item(47, 6).Value = 54  # F47: 53 -> 54

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 159  # F6: 157 -> 159
$ws.Cells.Item(19, 6).Value = 95  # F19: 94 -> 95
$ws.Cells.Item(23, 6).Value = 126  # F23: 124 -> 126
$ws.Cells.Item(26, 6).Value = 224  # F26: 223 -> 224

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 6).Value = 427  # F5: 426 -> 427
$ws.Cells.Item(7, 6).Value = 1513  # F7: 1512 -> 1513
$ws.Cells.Item(9, 6).Value = 443  # F9: 442 -> 443
$ws.Cells.Item(10, 6).Value = 2981  # F10: 2974 -> 2981
$ws.Cells.Item(13, 6).Value = 917  # F13: 913 -> 917
$ws.Cells.Item(14, 6).Value = 920  # F14: 914 -> 920
$ws.Cells.Item(15, 6).Value = 1410  # F15: 1406 -> 1410

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 427  # F2: 426 -> 427
$ws.Cells.Item(5, 6).Value = 442  # F5: 439 -> 442
$ws.Cells.Item(6, 6).Value = 443  # F6: 442 -> 443
$ws.Cells.Item(7, 6).Value = 2981  # F7: 2974 -> 2981
$ws.Cells.Item(9, 6).Value = 2895  # F9: 2891 -> 2895
$ws.Cells.Item(11, 6).Value = 755  # F11: 754 -> 755
$ws.Cells.Item(12, 6).Value = 7278  # F12: 7261 -> 7278
$ws.Cells.Item(13, 6).Value = 178  # F13: 177 -> 178
$ws.Cells.Item(16, 6).Value = 191  # F16: 190 -> 191
$ws.Cells.Item(17, 6).Value = 1417  # F17: 1412 -> 1417
$ws.Cells.Item(18, 6).Value = 920  # F18: 914 -> 920
$ws.Cells.Item(20, 6).Value = 127  # F20: 126 -> 127
$ws.Cells.Item(21, 6).Value = 3317  # F21: 3298 -> 3317
$ws.Cells.Item(23, 6).Value = 5771  # F23: 5764 -> 5771
$ws.Cells.Item(24, 6).Value = 564  # F24: 560 -> 564
$ws.Cells.Item(25, 6).Value = 928  # F25: 923 -> 928
$ws.Cells.Item(28, 6).Value = 5736  # F28: 5721 -> 5736
$ws.Cells.Item(29, 6).Value = 319  # F29: 318 -> 319
$ws.Cells.Item(30, 6).Value = 3948  # F30: 3941 -> 3948
$ws.Cells.Item(31, 6).Value = 659  # F31: 657 -> 659
$ws.Cells.Item(33, 6).Value = 1809  # F33: 1806 -> 1809
$ws.Cells.Item(35, 6).Value = 126  # F35: 124 -> 126
$ws.Cells.Item(36, 6).Value = 142  # F36: 141 -> 142
$ws.Cells.Item(37, 6).Value = 72  # F37: 66 -> 72
$ws.Cells.Item(38, 6).Value = 300  # F38: 298 -> 300
$ws.Cells.Item(39, 6).Value = 1097  # F39: 1096 -> 1097
$ws.Cells.Item(40, 6).Value = 1807  # F40: 1801 -> 1807
$ws.Cells.Item(41, 6).Value = 69  # F41: 68 -> 69
$ws.Cells.Item(42, 6).Value = 342  # F42: 337 -> 342
$ws.Cells.Item(44, 6).Value = 989  # F44: 985 -> 989
